# archive api: add role module
#
# Appends two new bulleted paragraphs right after the paragraph that
# contains "Модуль удаления записи" (the last item of the "Подсистема
# репозиториев" subsystem):
#
#   * Подсистема ролей                                            (ilvl 0)
#       * Модуль установки ограничения на выполнение метода по ролям (ilvl 1)
#
# Both paragraphs reuse the same pStyle/numPr/rFonts/sz formatting that is
# already used throughout the document's bullet list.

$d = $word.ActiveDocument

# Anchor on the paragraph that currently ends the document ("Модуль
# удаления записи") by locating its text rather than a hard-coded index.
$anchor = $d.Content
$null = $anchor.Find.Execute("Модуль удаления записи", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
$anchor.Expand(4)   # wdParagraph -> whole paragraph, including its mark
$anchor.Collapse(0) # collapse to the end of that paragraph

# --- New paragraph 1: "Подсистема ролей" (ilvl 0) -------------------------
$anchor.InsertParagraphAfter()
$para1 = $anchor.Paragraphs(1)

$frag1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t>Подсистема ролей</w:t></w:r></w:p>'
$para1.Range.InsertXML($frag1)

# --- New paragraph 2: "Модуль установки ограничения ..." (ilvl 1) --------
$r1 = $para1.Range
$r1.Collapse(0)
$r1.InsertParagraphAfter()
$para2 = $r1.Paragraphs(1)

$frag2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Модуль </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t>установки ограничения на выполнение метода по ролям</w:t></w:r></w:p>'
$para2.Range.InsertXML($frag2)

Write-Host "Inserted role module paragraphs."
